$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 (shifts K_CALCMETH..K_ZUORDN down by one)
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the K_BEWERTUNG entry
$ws.Cells.Item(6, 1).Value = "K_BEWERTUNG"
$ws.Cells.Item(6, 2).Value = "Grad der Effektivität"
$ws.Cells.Item(6, 3).Value = "Level of effectiveness"

# Copy formatting (style) from the row above so the new row matches existing data rows
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Widen column A per the diff (target stored width 33.54296875 chars;
# the nearest value reachable through the COM ColumnWidth pixel grid is used)
$ws.Columns.Item(1).ColumnWidth = 32.857142857142854
